$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.836.26"
$ws.Range("E2").Value = "  +0.23%  "

$ws.Range("D3").Value = "1.643.04"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("E4").Value = "  -0.75%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.27"
$ws.Range("E5").Value = "  -0.89%  "

$ws.Range("E6").Value = "  +0.94%  "

$ws.Range("E8").Value = "  +1.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0620"
$ws.Range("E9").Value = "  -0.39%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.82"
$ws.Range("E10").Value = "  +3.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("D12").Value = "1.872.92"
$ws.Range("E12").Value = "  +0.06%  "

$ws.Range("D13").Value = "1.657.20"
$ws.Range("E13").Value = "  +0.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.12"
$ws.Range("E14").Value = "  -0.31%  "

$ws.Range("E15").Value = "  +0.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.40"
$ws.Range("E16").Value = "  +2.67%  "

$ws.Range("D17").Value = "26.850.11"
$ws.Range("E17").Value = "  +0.37%  "

$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("E18").Value = "  +0.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.44"
$ws.Range("E19").Value = "  +3.02%  "

$ws.Range("E20").Value = "  -0.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.38"
$ws.Range("E21").Value = "  +1.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.61"
$ws.Range("E22").Value = "  +6.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.44"
$ws.Range("E23").Value = "  +5.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.16"
$ws.Range("E24").Value = "  -0.91%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.45"
$ws.Range("E25").Value = "  -1.44%  "

$ws.Range("E26").Value = "  -0.75%  "

$ws.Range("E27").Value = "  +4.21%  "

$ws.Range("E28").Value = "  +0.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.86"
$ws.Range("E29").Value = "  +1.88%  "

$ws.Range("E30").Value = "  +1.97%  "

$ws.Range("E31").Value = "  -0.67%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.37"
$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("E33").Value = "  +0.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.55"
$ws.Range("E34").Value = "  +1.81%  "

$ws.Range("E35").Value = "  -0.60%  "

$ws.Range("D36").Value = "1.242.20"
$ws.Range("E36").Value = "  -2.49%  "

$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.539"
$ws.Range("E38").Value = "  +2.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.832"
$ws.Range("E39").Value = "  +3.23%  "

$ws.Range("E40").Value = "  -0.77%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.805"
$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("E42").Value = "  +1.76%  "

$ws.Range("D43").Value = "1.785.89"
$ws.Range("E43").Value = "  +0.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "60.95"
$ws.Range("E45").Value = "  +1.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.55"
$ws.Range("E46").Value = "  +0.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.58"
$ws.Range("E47").Value = "  +0.12%  "

$ws.Range("E48").Value = "  +1.33%  "

$ws.Range("E49").Value = "  -1.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0972"
$ws.Range("E50").Value = "  +1.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.56"
$ws.Range("E51").Value = "  +0.50%  "
